$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$rows = @{
    12 = @{ H=0; I=0; J=0; K=0; L=0; N=$null }
    40 = @{ H=3265.5; I=2674.625; J=3480.3635; K=2674.625; L=3480.3635; M=-2499.625; N=-3830.3635 }
    74 = @{ H=3461.9048; I=3200; J=3489.4736; K=3200; L=3489.4736; M=-2264; N=-5361.473599999999 }
    77 = @{ H=3461.9048; I=3200; J=3489.4736; K=16000; L=17447.368; M=-11320; N=-26807.368 }
    80 = @{ H=633.26666; I=700; J=588.7778; K=2100; L=1766.3334; M=-1102; N=-3762.3334 }
    83 = @{ H=633.26666; I=700; J=588.7778; K=6300; L=5299.000199999999; M=-1308; N=-15283.0002 }
    101 = @{ H=7333.727; I=511.33334; J=12056.923; K=1534.00002; L=36170.769; M=87.99998000000005; N=-39414.769 }
}
foreach ($r in $rows.Keys) {
    $cols = $rows[$r]
    foreach ($c in $cols.Keys) {
        $addr = "$c$r"
        $ws.Range($addr).Value = $cols[$c]
    }
}

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$rows = @{
    32 = @{ H=12831243; I=16952430; J=33872.58; K=16952430; L=33872.58; M=-16952143; N=-34446.58 }
}
foreach ($r in $rows.Keys) {
    $cols = $rows[$r]
    foreach ($c in $cols.Keys) {
        $addr = "$c$r"
        $ws.Range($addr).Value = $cols[$c]
    }
}

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$rows = @{
    20 = @{ H=4200.778; I=4002.6667; J=4299.8335; K=4002.6667; L=4299.8335; M=-3755.6667; N=-4793.8335 }
    28 = @{ H=30000; I=0; J=30000; K=0; L=30000; N=-30588 }
    86 = @{ H=1601.76; I=1755; J=1406.7273; K=1755; L=1406.7273; M=-632; N=-3652.7273 }
    89 = @{ H=1601.76; I=1755; J=1406.7273; K=8775; L=7033.636500000001; M=-3159; N=-18265.6365 }
    107 = @{ H=3039.682; I=3043.6; J=3000.5; K=3043.6; L=3000.5; M=-1123.6; N=-6840.5 }
}
foreach ($r in $rows.Keys) {
    $cols = $rows[$r]
    foreach ($c in $cols.Keys) {
        $addr = "$c$r"
        $ws.Range($addr).Value = $cols[$c]
    }
}

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$rows = @{
    22 = @{ H=267.8; I=159; J=567; K=159; L=567; M=191; N=-1267 }
    31 = @{ H=3033.2812; I=1823.3478; J=6125.3335; K=1823.3478; L=6125.3335; M=-1528.3478; N=-6715.3335 }
    34 = @{ H=3033.2812; I=1823.3478; J=6125.3335; K=1823.3478; L=6125.3335; M=-1621.3478; N=-6529.3335 }
    58 = @{ H=3145.85; I=1672; J=4351.727; K=1672; L=4351.727; M=-1469; N=-4757.727 }
    76 = @{ H=5000; I=5000; J=0; K=5000; L=0; M=-4685 }
    79 = @{ H=5000; I=5000; J=0; K=5000; L=0; M=-3908 }
    86 = @{ H=3799.4; I=3003.5; J=4330; K=3003.5; L=4330; M=-1880.5; N=-6576 }
    89 = @{ H=3799.4; I=3003.5; J=4330; K=15017.5; L=21650; M=-9401.5; N=-32882 }
    136 = @{ H=3145.85; I=1672; J=4351.727; K=5016; L=13055.181; M=-2466; N=-18155.181 }
}
foreach ($r in $rows.Keys) {
    $cols = $rows[$r]
    foreach ($c in $cols.Keys) {
        $addr = "$c$r"
        $ws.Range($addr).Value = $cols[$c]
    }
}

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$rows = @{
    122 = @{ H=685.09375; J=889.8; L=8008.2; N=-12908.2 }
    131 = @{ H=824.7778; I=655.4; J=863.2727; K=1966.2; L=2589.8181; M=3073.8; N=-12669.8181 }
}
foreach ($r in $rows.Keys) {
    $cols = $rows[$r]
    foreach ($c in $cols.Keys) {
        $addr = "$c$r"
        $ws.Range($addr).Value = $cols[$c]
    }
}

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$rows = @{
    2 = @{ H=61.9; I=64.333336; J=40; K=64.333336; L=40; M=48.666664; N=-266 }
    97 = @{ H=1129.9565; I=696; J=1943.625; K=696; L=1943.625; M=-200; N=-2935.625 }
    107 = @{ H=155.33333; J=70; L=70; N=-3910 }
}
foreach ($r in $rows.Keys) {
    $cols = $rows[$r]
    foreach ($c in $cols.Keys) {
        $addr = "$c$r"
        $ws.Range($addr).Value = $cols[$c]
    }
}

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$rows = @{
    46 = @{ H=896.6667; I=0; J=896.6667; K=0; L=896.6667; M=$null; N=-1272.6667 }
    99 = @{ H=17666.666; I=16000; J=21000; K=16000; L=21000; M=-13005; N=-26990 }
}
foreach ($r in $rows.Keys) {
    $cols = $rows[$r]
    foreach ($c in $cols.Keys) {
        $addr = "$c$r"
        $ws.Range($addr).Value = $cols[$c]
    }
}

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$rows = @{
    123 = @{ H=18286.666; J=18286.666; L=18286.666; N=-28086.666 }
}
foreach ($r in $rows.Keys) {
    $cols = $rows[$r]
    foreach ($c in $cols.Keys) {
        $addr = "$c$r"
        $ws.Range($addr).Value = $cols[$c]
    }
}
